# Trade #24 closed at 2026-02-17 20:07:22 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.82
$summary.Range("B4").Value = -0.19
$summary.Range("B6").Value = 24
$summary.Range("B8").Value = 11
$summary.Range("B9").Value = 50

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.81999999999999
$status.Range("D5").Value = 24
$status.Range("E5").Value = -0.19
$status.Range("F5").Value = -0.18
$status.Range("G5").Value = 50

# ---------------------------------------------------------------------------
# New trade row (#24) appended to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------------
$newRow = @(24, "2026-02-17", "20:07:15", "MarketMaking", "DOWN", 0.05, 0.04, "CLOSED", -20, -0.01, 99.81999999999999, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 25
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($row, $col)
        if ($col -eq 2) {
            # Column B holds an ISO-looking date string ("2026-02-17") that
            # must stay plain text, not get auto-converted into a date
            # serial number. Force the cell to Text first, then clear the
            # formatting back off so it ends up styled like its neighbours.
            $cell.NumberFormat = "@"
            $cell.Value = $newRow[$i]
            $cell.ClearFormats()
        } else {
            $cell.Value = $newRow[$i]
        }
    }
}
